$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet updates ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: $newVersion"

$aboutWs.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Wangjiayu Coal Mine, China, M2184, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 7; $r++) {
    $dataWs.Range("S$r").Value = $newVersion
}
